# Update "想去人数" (column F) figures on the sheets that list event data.
# Both the "展览" sheet and the "全部类型" sheet carry the same rows, so the
# same updates are applied to each.

$wb = $excel.ActiveWorkbook

# Row number (1-based, matching worksheet row) -> new value for column F.
$updates = @{
    3  = 1075
    8  = 11144
    9  = 4276
    14 = 1067
    15 = 93
    18 = 482
    19 = 11220
    20 = 11064
    22 = 36
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
